$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the JSON output paths: "New" -> "Json"
$ws.Range("B2").Value = "D:\PythonProjects\TestExcelToJson\Json\1.json"
$ws.Range("B3").Value = "D:\PythonProjects\TestExcelToJson\Json\2.json"
$ws.Range("B4").Value = "D:\PythonProjects\TestExcelToJson\Json\3.json"

# Move the active selection from B8 to C6
$ws.Range("C6").Select()
